# Notes taken on 6/20/22
$wb = $excel.ActiveWorkbook

# --- "stocks" sheet: reorder columns B:E from (tickers, price, pe, eps)
#     to (eps, pe, price, tickers) ---
$ws1 = $wb.Worksheets.Item("stocks")

$stockRows = @{}
for ($r = 2; $r -le 4; $r++) {
    $stockRows[$r] = @(
        $ws1.Cells.Item($r, 2).Value2,  # B - ticker
        $ws1.Cells.Item($r, 3).Value2,  # C - price
        $ws1.Cells.Item($r, 4).Value2,  # D - pe
        $ws1.Cells.Item($r, 5).Value2   # E - eps
    )
}

$ws1.Range("B1").Value = "eps"
$ws1.Range("C1").Value = "pe"
$ws1.Range("D1").Value = "price"
$ws1.Range("E1").Value = "tickers"

for ($r = 2; $r -le 4; $r++) {
    $vals = $stockRows[$r]
    $ws1.Cells.Item($r, 2).Value = $vals[3]  # B = eps
    $ws1.Cells.Item($r, 3).Value = $vals[2]  # C = pe
    $ws1.Cells.Item($r, 4).Value = $vals[1]  # D = price
    $ws1.Cells.Item($r, 5).Value = $vals[0]  # E = tickers
}

# --- "weather" sheet: swap the "temperature" (C) and "event" (D) columns ---
$ws2 = $wb.Worksheets.Item("weather")

$weatherRows = @{}
for ($r = 2; $r -le 4; $r++) {
    $weatherRows[$r] = @(
        $ws2.Cells.Item($r, 3).Value2,  # C - temperature
        $ws2.Cells.Item($r, 4).Value2   # D - event
    )
}

$ws2.Range("C1").Value = "event"
$ws2.Range("D1").Value = "temperature"

for ($r = 2; $r -le 4; $r++) {
    $vals = $weatherRows[$r]
    $ws2.Cells.Item($r, 3).Value = $vals[1]  # C = event
    $ws2.Cells.Item($r, 4).Value = $vals[0]  # D = temperature
}
